{"js": "const body = context.document.body;\n\n// ---- Change 1: bold \"REDO APPLY\" inside the PHYSICAL STANDBY paragraph ----\nconst redoApply = body.search(\"REDO APPLY\", { matchCase: true });\nawait context.sync();\nredoApply.items[0].font.bold = true;\nawait context.sync();\n\n// ---- Change 2: \"Where we can do read,write operations\" ----\n// becomes \"Where we can do **read, write** operations\" (comma gets a space,\n// and \"read, write\" becomes bold). We also cleanly separate \"Where\" onto its\n// own run (matching the restructured runs in the target document) by\n// round-tripping its highlight color, which forces a run split without\n// altering the visible formatting.\nconst whereRange = body.search(\"Where\", { matchCase: true });\nawait context.sync();\nwhereRange.items[0].font.highlightColor = \"#00FF00\";\nawait context.sync();\nwhereRange.items[0].font.highlightColor = \"#FFFF00\";\nawait context.sync();\n\n// Fix comma spacing: \"read,write\" -> \"read, write\"\nconst readWrite = body.search(\"read,write\", { matchCase: true });\nawait context.sync();\nreadWrite.items[0].insertText(\"read, write\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Bold \"read, write\"\nconst readWriteBold = body.search(\"read, write\", { matchCase: true });\nawait context.sync();\nreadWriteBold.items[0].font.bold = true;\nawait context.sync();\n\n// ---- Change 3: bold \"BOTH READ/WRITE\" inside the SNAPSHOT STANDBY paragraph ----\nconst bothReadWrite = body.search(\"BOTH READ/WRITE\", { matchCase: true });\nawait context.sync();\nbothReadWrite.items[0].font.bold = true;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---- Change 1: bold \"REDO APPLY\" inside the PHYSICAL STANDBY paragraph ----\n$rng1 = $d.Content\n$rng1.Find.Text = \"REDO APPLY\"\n$rng1.Find.MatchCase = $true\n$rng1.Find.Execute() | Out-Null\n$rng1.Bold = 1\n\n# ---- Change 2: \"Where we can do read,write operations\" ----\n# becomes \"Where we can do **read, write** operations\" (comma gets a space,\n# and \"read, write\" becomes bold). First cleanly separate \"Where\" onto its\n# own run (matching the restructured runs in the target document) by\n# round-tripping its font size, which forces a run split without altering\n# the visible formatting.\n$rng2 = $d.Content\n$rng2.Find.Text = \"Where\"\n$rng2.Find.MatchCase = $true\n$rng2.Find.Execute() | Out-Null\n$rng2.Font.Size = 11\n$rng2.Font.Size = 12\n\n# Fix comma spacing: \"read,write\" -> \"read, write\"\n$rng3 = $d.Content\n$rng3.Find.Text = \"read,write\"\n$rng3.Find.MatchCase = $true\n$rng3.Find.Execute() | Out-Null\n$rng3.Text = \"read, write\"\n\n# Bold \"read, write\"\n$rng4 = $d.Content\n$rng4.Find.Text = \"read, write\"\n$rng4.Find.MatchCase = $true\n$rng4.Find.Execute() | Out-Null\n$rng4.Bold = 1\n\n# ---- Change 3: bold \"BOTH READ/WRITE\" inside the SNAPSHOT STANDBY paragraph ----\n$rng5 = $d.Content\n$rng5.Find.Text = \"BOTH READ/WRITE\"\n$rng5.Find.MatchCase = $true\n$rng5.Find.Execute() | Out-Null\n$rng5.Bold = 1\n"}
